$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.291273832321167
$ws.Range("B1").Value = 2.924385070800781
$ws.Range("C1").Value = 5.261811256408691
$ws.Range("D1").Value = 1.85080349445343
$ws.Range("E1").Value = 1.01287305355072
